$wb = $excel.ActiveWorkbook

# --- connectivity_score (sheet4): selection moves from H11 to B1 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B1").Select() | Out-Null

# --- add new eia_area sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "eia_area"

# --- populate segment_lookup table ---
$ws5.Range("A2").Value = 'Rock Creek - Upper - MS4 - Portal Branch'
$ws5.Range("B2").Value = 'TPO01'
$ws5.Range("A3").Value = 'Rock Creek - Lower - MS4 - Dumbarton Oaks'
$ws5.Range("B3").Value = 'TDO01'
$ws5.Range("A4").Value = 'Rock Creek - Lower - MS4 - Klingle Valley Run'
$ws5.Range("B4").Value = 'TKV01'
$ws5.Range("A5").Value = 'Rock Creek - Lower - MS4 - Melvin Hazen Valley Branch'
$ws5.Range("B5").Value = 'TMH01'
$ws5.Range("A6").Value = 'Rock Creek - Lower - MS4 - Normanstone Creek'
$ws5.Range("B6").Value = 'TNS01'
$ws5.Range("A7").Value = 'Potomac River - Upper - MS4 - Battery Kemble Creek'
$ws5.Range("B7").Value = 'TBK01'
$ws5.Range("A8").Value = 'Anacostia River - Upper - MS4 - Fort Davis Tributary'
$ws5.Range("B8").Value = 'TFD01'
$ws5.Range("A9").Value = 'Anacostia River - Upper - MS4 - Texas Avenue Tributary'
$ws5.Range("B9").Value = 'TTX27'
$ws5.Range("A10").Value = 'Anacostia River - Lower - MS4 - Fort Stanton Tributary'
$ws5.Range("B10").Value = 'TFS01'
$ws5.Range("A11").Value = 'Potomac River - Upper - MS4 - Dalecarlia Tributary'
$ws5.Range("B11").Value = 'TDA01'
$ws5.Range("A12").Value = 'Anacostia River - Upper - MS4 - Pope Branch'
$ws5.Range("B12").Value = 'TPB01'
$ws5.Range("A13").Value = 'Anacostia River - Upper - MS4 - Fort Chaplin Tributary'
$ws5.Range("B13").Value = 'TFC01'
$ws5.Range("A14").Value = 'Anacostia River - Upper - MS4 - Fort Dupont Tributary'
$ws5.Range("B14").Value = 'TDU01'
$ws5.Range("A15").Value = 'Rock Creek - Upper - MS4 - Soapstone Creek'
$ws5.Range("B15").Value = 'TSO01'
$ws5.Range("A16").Value = 'Anacostia River - Upper - MS4 - Nash Run'
$ws5.Range("B16").Value = 'TNA01'
$ws5.Range("A17").Value = 'Rock Creek - Upper - MS4 - Luzon Branch'
$ws5.Range("B17").Value = 'TLU01'
$ws5.Range("A18").Value = 'Rock Creek - Upper - MS4 - Pinehurst Branch'
$ws5.Range("B18").Value = 'TPI01'
$ws5.Range("A19").Value = 'Rock Creek - Upper - MS4 - Fenwick Branch'
$ws5.Range("B19").Value = 'TFE01'
$ws5.Range("A20").Value = 'Anacostia River - Upper - MS4 - Hickey Run'
$ws5.Range("B20").Value = 'THR01'
$ws5.Range("A21").Value = 'Potomac River - Upper - MS4 - Foundry Branch'
$ws5.Range("B21").Value = 'TFB01'
$ws5.Range("A22").Value = 'Rock Creek - Upper - MS4 - Broad Branch'
$ws5.Range("B22").Value = 'TBR01'
$ws5.Range("A23").Value = 'Rock Creek - Lower - CSS - Rock Creek'
$ws5.Range("B23").Value = 'RCR01'
$ws5.Range("B24").Value = 'TWB06'
$ws5.Range("A24").Value = 'Anacostia River - Upper - MS4 - Watts Branch - Upper'
$ws5.Range("B25").Value = 'TWB05'
$ws5.Range("A25").Value = 'Anacostia River - Upper - MS4 - Watts Branch - Lower'
$ws5.Range("A26").Value = 'Rock Creek - Lower - MS4 - Piney Branch'
$ws5.Range("B26").Value = 'TPY01'
$ws5.Range("A27").Value = 'Potomac River - Lower - MS4 - Oxon Run'
$ws5.Range("B27").Value = 'TOR01'
$ws5.Range("B28").Value = 'RCR05'
$ws5.Range("A28").Value = 'Rock Creek - Upper - MS4 - Rock Creek'
$ws5.Range("A29").Value = 'Rock Creek - Lower - MS4 - Rock Creek'
$ws5.Range("B29").Value = 'RCR09'
$ws5.Range("A1").Value = 'segment_lookup'
$ws5.Range("B1").Value = 'location_id'

# --- set column widths (approx bestFit, quantized to nearest achievable step) ---
$ws5.Columns.Item(1).ColumnWidth = 47.5
$ws5.Columns.Item(2).ColumnWidth = 9.5

# --- final selection/activation on eia_area ---
$ws5.Range("G14").Select() | Out-Null
